$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update row 2 values
$ws.Range("G2").Value = 0.05014
$ws.Range("H2").Value = 0.10028
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.2696695
$ws.Range("N2").Value = 2.539339
$ws.Range("O2").Value = 0.2673143181973693
$ws.Range("P2").Value = 0.2673143181973693
$ws.Range("Q2").Value = 0.06366122872999999
$ws.Range("R2").Value = 0.25464491492
$ws.Range("S2").Value = 0.2673143181973693
$ws.Range("T2").Value = 0.2673143181973693

# Update row 3 values
$ws.Range("G3").Value = 0.05014
$ws.Range("H3").Value = 0.10028
$ws.Range("M3").Value = 3.4800555
$ws.Range("N3").Value = 6.960110999999999
$ws.Range("O3").Value = 0.7326856818026307
$ws.Range("P3").Value = 0.7326856818026307
$ws.Range("Q3").Value = 0.17448998277
$ws.Range("R3").Value = 0.6979599310799999
$ws.Range("S3").Value = 0.7326856818026307
$ws.Range("T3").Value = 0.7326856818026307

# Delete entire row 4 (shifts rows up, removes the row completely)
$ws.Rows.Item(4).Delete()
